# Update column F (dSF) values for rows 2-10 per repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = -9
$ws.Range("F9").Value = 13
$ws.Range("F10").Value = -3
